# Update cryptocurrency price/volume table (columns D and E) to reflect
# the latest scrape, per the GitHub Actions scheduled refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: column D holds free-form price strings (type text in the source
# data, e.g. "34.754.46" or "228.43"). Some of these parse as plain
# numbers, so a leading apostrophe is used to force text entry and keep
# them stored as strings (matching the original file), exactly like
# typing '228.43 into a General-formatted cell in Excel.

$ws.Range("D2").Value = "'34.754.46"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("D3").Value = "'1.818.21"
$ws.Range("E3").Value = "  +1.50%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'228.43"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("E6").Value = "  +1.60%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "'35.17"
$ws.Range("E8").Value = "  +8.99%  "
$ws.Range("E9").Value = "  +1.90%  "
$ws.Range("D10").Value = "'0.0698"
$ws.Range("E10").Value = "  +1.27%  "
$ws.Range("D11").Value = "'0.0952"
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("D12").Value = "'2.080.77"
$ws.Range("E12").Value = "  +1.52%  "
$ws.Range("D13").Value = "'11.41"
$ws.Range("E13").Value = "  +4.17%  "
$ws.Range("D14").Value = "'1.818.48"
$ws.Range("E14").Value = "  +1.44%  "
$ws.Range("D15").Value = "'0.645"
$ws.Range("E15").Value = "  +2.86%  "
$ws.Range("D16").Value = "'34.712.02"
$ws.Range("E16").Value = "  +1.17%  "
$ws.Range("D17").Value = "'4.35"
$ws.Range("E17").Value = "  +3.40%  "
$ws.Range("E18").Value = "  +1.79%  "
$ws.Range("D19").Value = "'249.09"
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("D20").Value = "'0.0₃0804"
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("D21").Value = "'11.59"
$ws.Range("E21").Value = "  +5.84%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "'4.22"
$ws.Range("E23").Value = "  +1.09%  "
$ws.Range("D24").Value = "'172.71"
$ws.Range("E24").Value = "  +6.29%  "
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("D26").Value = "'7.48"
$ws.Range("E26").Value = "  +4.39%  "
$ws.Range("D27").Value = "'16.86"
$ws.Range("E27").Value = "  +3.04%  "
$ws.Range("E28").Value = "  +1.49%  "
$ws.Range("D30").Value = "'4.01"
$ws.Range("E30").Value = "  +3.33%  "
$ws.Range("D31").Value = "'0.0534"
$ws.Range("E31").Value = "  +2.54%  "
$ws.Range("D32").Value = "'3.87"
$ws.Range("E32").Value = "  +2.61%  "
$ws.Range("E33").Value = "  +1.62%  "
$ws.Range("D34").Value = "'1.87"
$ws.Range("E34").Value = "  +3.16%  "
$ws.Range("D35").Value = "'2.67"
$ws.Range("E35").Value = "  +1.79%  "
$ws.Range("D36").Value = "'1.425.12"
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("D37").Value = "'0.681"
$ws.Range("E37").Value = "  +3.10%  "
$ws.Range("E38").Value = "  +1.84%  "
$ws.Range("D39").Value = "'86.55"
$ws.Range("E39").Value = "  +5.32%  "
$ws.Range("E40").Value = "  +0.76%  "
$ws.Range("E41").Value = "  +4.14%  "
$ws.Range("D42").Value = "'0.967"
$ws.Range("E42").Value = "  +4.26%  "
$ws.Range("D43").Value = "'2.40"
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("D44").Value = "'13.84"
$ws.Range("E44").Value = "  -1.76%  "
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("E46").Value = "  +3.01%  "
$ws.Range("D47").Value = "'6.15"
$ws.Range("E47").Value = "  +1.74%  "
$ws.Range("D48").Value = "'1.981.04"
$ws.Range("D49").Value = "'106.61"
$ws.Range("E49").Value = "  +1.14%  "
$ws.Range("D50").Value = "'0.0₆0132"
$ws.Range("E50").Value = "  +0.89%  "
$ws.Range("E51").Value = "  -0.07%  "
